$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Corrected ICDC "FilesTab" Cypher query (B4): drop the File Type and Breed
# columns from the RETURN clause.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bulldog']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B4").Value = $newQuery

# Column B is narrower now that the query text is shorter.
$ws.Columns.Item(2).ColumnWidth = 67.5

# Row 4 shrinks to match the new (shorter) wrapped query text.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection / scroll position moved down to the corrected row.
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
